# "updated UI for batch operation"
# - Rename existing sheet "Sheet1" -> "v1"
# - Add a new sheet "v2" after it, with the same style of milestone table
#   (columns A/B/F) but populated with the new "small-cell repowering"
#   batch-operation milestone list instead of the old one.
# - v2 becomes the active / selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- rename + add the new sheet right after v1 ---------------------------
$ws1.Name = "v1"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "v2"

# --- milestone names for the new batch-operation workflow ----------------
$names = @(
    "repowering_up",
    "inter_transmission_merge",
    "rfi",
    "rfnsa_check",
    "acma_check",
    "bbu_status_check",
    "shutdown_cr",
    "ssv_pre_cutover",
    "naming_convention",
    "bbu_cutover_cr",
    "bbu_cutover",
    "site_list_check",
    "overlap_simulation",
    "pci_conflict",
    "rfnsa_update",
    "acma_update",
    "prs_cellgroup",
    "rf_script",
    "activation_cr",
    "emeg_check",
    "cell_activation",
    "ssv_post_cutover",
    "service_notification",
    "isn_upload",
    "dsa7_upload"
)

# rows 2..26 hold the 25 milestone names
$firstRow = 2
$lastRow = $firstRow + $names.Count - 1

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $firstRow + $i
    $ws2.Cells.Item($row, 1).Value = $names[$i]
}

# column B: per-row "'name', " helper formula -- first two rows get their
# own (non-shared) formula, the remaining block B4:B26 is filled as one go
# so the engine groups it into a single shared formula (matches how the
# original v1 sheet was built).
$ws2.Range("B2").Formula = "=""'""&A2&""', """
$ws2.Range("B3").Formula = "=""'""&A3&""', """
$ws2.Range("B4:B26").Formula = "=""'""&A4&""', """

# --- highlight the two "cell_activation" / "ssv_post_cutover" rows -------
# (same yellow fill used elsewhere in the workbook for a highlighted group)
$ws2.Range("A22:A23").Interior.Color = 65535

# --- column F: concatenated "quote-comma" strings for each logical group -
$ws2.Range("F2").Formula = "=B2&B3&B4&B5&B6&B7&B8&B9&B10&B11&B12&B13&B14&B15&B16&B17&B18&B19&B20&B21"
$ws2.Range("F22").Formula = "=B22&B23"
$ws2.Range("F24").Formula = "=B24&B25&B26"

# recalc so F2/F22/F24 have their text available to copy as static notes
$wb.Application.Calculate()

# F3 / F23 / F25 are the typed (quote-prefixed) continuation notes --
# highlighted with the same yellow fill as the F-column "group" markers
# elsewhere in the workbook, and typed as literal text starting with a
# leading apostrophe (quote-prefixed text), matching the style used for
# F4 / F32 on v1.
$ws2.Range("F3").Interior.Color = 65535
$ws2.Range("F3").Value = "'" + $ws2.Range("F2").Value()

$ws2.Range("F23").Interior.Color = 65535
$ws2.Range("F23").Value = "'" + $ws2.Range("F22").Value()

$ws2.Range("F25").Interior.Color = 65535
$ws2.Range("F25").Value = "'" + $ws2.Range("F24").Value()

# --- column widths (closest achievable to the authored widths) -----------
$ws2.Columns.Item(1).ColumnWidth = 26
$ws2.Columns.Item(2).ColumnWidth = 10.75

# --- selection / active-sheet state ---------------------------------------
$ws1.Activate()
$ws1.Range("B3").Select()

$ws2.Activate()
$ws2.Range("U16").Select()
